$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume/1h change (column E) values
# as refreshed by the GitHub Actions data pull on 2023-01-15.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '301.85'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.85%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.46'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.84%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.145'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-3.25%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07388'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.12%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.832'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '24.94%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.860'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.61%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.752'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-2.82%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9289'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.75%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1710'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.36%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07247'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-6.10%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08099'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.39%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03028'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-1.21%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09955'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.69%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001491'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.36%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006086'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-5.85%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.27%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.227'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.08%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3259'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.31%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-2.12%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.579'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.39%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04649'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.67%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1579'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-2.86%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001216'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.17%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004778'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '7.99%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.44%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '7.42%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01726'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-1.50%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04524'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.44%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007083'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.84%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1347'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.28%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002143'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.03%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01046'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-17.68%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006219'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '0.82%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-21.49%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.844'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '160.42%'
